# Update cell values in the Leve profit tables across multiple sheets
# (data refreshed from a scheduled market-board pull).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 635.3570999999999
$ws.Range("I12").Value = 426.125
$ws.Range("K12").Value = 426.125
$ws.Range("M12").Value = -256.125
$ws.Range("H19").Value = 39395692
$ws.Range("I19").Value = 66668068
$ws.Range("J19").Value = 16668716
$ws.Range("K19").Value = 66668068
$ws.Range("L19").Value = 16668716
$ws.Range("M19").Value = -66667893
$ws.Range("N19").Value = -16669066
$ws.Range("H101").Value = 871.7692
$ws.Range("J101").Value = 1347
$ws.Range("L101").Value = 4041
$ws.Range("N101").Value = -7285
$ws.Range("H103").Value = 62502256
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 62502256
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 187506768
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -187507940
$ws.Range("H132").Value = 2639.3572
$ws.Range("I132").Value = 2663.2222
$ws.Range("K132").Value = 7989.6666
$ws.Range("M132").Value = -5459.6666
$ws.Range("H135").Value = 2584.25
$ws.Range("I135").Value = 724.4
$ws.Range("K135").Value = 6519.599999999999
$ws.Range("M135").Value = -3984.599999999999
$ws.Range("H141").Value = 4488
$ws.Range("I141").Value = 4488
$ws.Range("K141").Value = 13464
$ws.Range("M141").Value = -8284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3256.9167
$ws.Range("I122").Value = 2521.4075
$ws.Range("J122").Value = 5463.4443
$ws.Range("K122").Value = 7564.2225
$ws.Range("L122").Value = 16390.3329
$ws.Range("M122").Value = -5114.2225
$ws.Range("N122").Value = -21290.3329
$ws.Range("H132").Value = 3759.5945
$ws.Range("I132").Value = 3517.3142
$ws.Range("K132").Value = 10551.9426
$ws.Range("M132").Value = -8021.942599999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3001.8
$ws.Range("I99").Value = 1010
$ws.Range("K99").Value = 1010
$ws.Range("M99").Value = 488
$ws.Range("H134").Value = 1891.1818
$ws.Range("I134").Value = 1891.1818
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5673.5454
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3138.5454
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 501.8
$ws.Range("I7").Value = 527.5
$ws.Range("K7").Value = 527.5
$ws.Range("M7").Value = -414.5
$ws.Range("H31").Value = 16396183
$ws.Range("I31").Value = 20835674
$ws.Range("K31").Value = 20835674
$ws.Range("M31").Value = -20835379
$ws.Range("H34").Value = 16396183
$ws.Range("I34").Value = 20835674
$ws.Range("K34").Value = 20835674
$ws.Range("M34").Value = -20835472
$ws.Range("H99").Value = 17152.541
$ws.Range("I99").Value = 10739.583
$ws.Range("J99").Value = 23565.5
$ws.Range("K99").Value = 10739.583
$ws.Range("L99").Value = 23565.5
$ws.Range("M99").Value = -9241.583000000001
$ws.Range("N99").Value = -26561.5
$ws.Range("H107").Value = 1191.125
$ws.Range("I107").Value = 366.0625
$ws.Range("K107").Value = 366.0625
$ws.Range("M107").Value = 1553.9375
$ws.Range("H122").Value = 3595.1667
$ws.Range("I122").Value = 1905.5
$ws.Range("K122").Value = 5716.5
$ws.Range("M122").Value = -3266.5
$ws.Range("H126").Value = 17152.541
$ws.Range("I126").Value = 10739.583
$ws.Range("J126").Value = 23565.5
$ws.Range("K126").Value = 32218.749
$ws.Range("L126").Value = 70696.5
$ws.Range("M126").Value = -29748.749
$ws.Range("N126").Value = -75636.5
$ws.Range("H132").Value = 2164.125
$ws.Range("I132").Value = 1936
$ws.Range("K132").Value = 5808
$ws.Range("M132").Value = -3278
$ws.Range("H134").Value = 2187.3076
$ws.Range("I134").Value = 2173.75
$ws.Range("J134").Value = 2350
$ws.Range("K134").Value = 6521.25
$ws.Range("L134").Value = 7050
$ws.Range("M134").Value = -3986.25
$ws.Range("N134").Value = -12120
$ws.Range("H141").Value = 554833.2
$ws.Range("J141").Value = 759999.75
$ws.Range("L141").Value = 759999.75
$ws.Range("N141").Value = -770359.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 3796663.5
$ws.Range("I107").Value = 2906.4
$ws.Range("J107").Value = 4795020.5
$ws.Range("K107").Value = 8719.200000000001
$ws.Range("L107").Value = 14385061.5
$ws.Range("M107").Value = -6799.200000000001
$ws.Range("N107").Value = -14388901.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2139.913
$ws.Range("I102").Value = 1792.1666
$ws.Range("J102").Value = 3391.8
$ws.Range("K102").Value = 1792.1666
$ws.Range("L102").Value = 3391.8
$ws.Range("M102").Value = -170.1666
$ws.Range("N102").Value = -6635.8
$ws.Range("H113").Value = 2583
$ws.Range("I113").Value = 2299.6667
$ws.Range("K113").Value = 2299.6667
$ws.Range("M113").Value = -129.6667000000002
$ws.Range("H132").Value = 6044.816
$ws.Range("I132").Value = 6046.36
$ws.Range("J132").Value = 6041.846
$ws.Range("K132").Value = 18139.08
$ws.Range("L132").Value = 18125.538
$ws.Range("M132").Value = -15609.08
$ws.Range("N132").Value = -23185.538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8808.5
$ws.Range("I7").Value = 8808.5
$ws.Range("K7").Value = 8808.5
$ws.Range("M7").Value = -8696.5
$ws.Range("H40").Value = 5368.28
$ws.Range("I40").Value = 4891.273
$ws.Range("K40").Value = 4891.273
$ws.Range("M40").Value = -4755.273
$ws.Range("H93").Value = 2647234
$ws.Range("I93").Value = 521.8182
$ws.Range("K93").Value = 521.8182
$ws.Range("M93").Value = 726.1818
$ws.Range("H122").Value = 3446.8628
$ws.Range("I122").Value = 3434.9302
$ws.Range("J122").Value = 3511
$ws.Range("K122").Value = 10304.7906
$ws.Range("L122").Value = 10533
$ws.Range("M122").Value = -7854.7906
$ws.Range("N122").Value = -15433
$ws.Range("H126").Value = 8808.5
$ws.Range("I126").Value = 8808.5
$ws.Range("K126").Value = 26425.5
$ws.Range("M126").Value = -23955.5
$ws.Range("H132").Value = 3527.743
$ws.Range("I132").Value = 2709.2
$ws.Range("J132").Value = 5574.1
$ws.Range("K132").Value = 8127.599999999999
$ws.Range("L132").Value = 16722.3
$ws.Range("M132").Value = -5597.599999999999
$ws.Range("N132").Value = -21782.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 130429
$ws.Range("J46").Value = 130429
$ws.Range("L46").Value = 130429
$ws.Range("N46").Value = -130891
$ws.Range("H122").Value = 3506
$ws.Range("I122").Value = 2916.818
$ws.Range("J122").Value = 5666.3335
$ws.Range("K122").Value = 8750.454000000002
$ws.Range("L122").Value = 16999.0005
$ws.Range("M122").Value = -6300.454000000002
$ws.Range("N122").Value = -21899.0005
$ws.Range("H132").Value = 1468.0186
$ws.Range("I132").Value = 968.3555
$ws.Range("J132").Value = 3966.3333
$ws.Range("K132").Value = 2905.0665
$ws.Range("L132").Value = 11898.9999
$ws.Range("M132").Value = -375.0664999999999
$ws.Range("N132").Value = -16958.9999
$ws.Range("H134").Value = 130429
$ws.Range("J134").Value = 130429
$ws.Range("L134").Value = 391287
$ws.Range("N134").Value = -396357
